$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added to the daily series. It belongs
# chronologically between the existing row 29 (2021-08-24) and the old
# row 30 (2021-08-20), so insert a fresh row at position 30 and shift
# everything from the old row 30 onward down by one (rows 30-37 -> 31-38).
$ws.Rows.Item(30).Insert()

# Populate the newly inserted row 30 with the new record.
$ws.Cells.Item(30, 1).Value = 9
$ws.Cells.Item(30, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(30, 3).Value = "Metropolitana"
$ws.Cells.Item(30, 4).Value = 44726
$ws.Cells.Item(30, 5).Value = 13
$ws.Cells.Item(30, 6).Value = 100112035
$ws.Cells.Item(30, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(30, 8).Value = "Sin especificar"
$ws.Cells.Item(30, 9).Value = "Primera"
$ws.Cells.Item(30, 10).Value = 28
$ws.Cells.Item(30, 11).Value = 24000
$ws.Cells.Item(30, 12).Value = 24000
$ws.Cells.Item(30, 13).Value = 24000
$ws.Cells.Item(30, 14).Value = "$/malla 15 kilos"
$ws.Cells.Item(30, 15).Value = "Hijuelas"
$ws.Cells.Item(30, 16).Value = 1600
$ws.Cells.Item(30, 17).Value = 15
$ws.Cells.Item(30, 18).Value = "Hortaliza"
